# edit.ps1
# Applies the changes described in the commit:
#   "Added OOI barcodes, NUTNR Cal params, FLORTO, corrected reference Desig."
#
# Summary of changes (from the OOXML diff):
#  - Sheet "Moorings": drop the stale cell selection (A2) left in the sheetView.
#  - Sheet "Asset_Cal_Info": the active-cell selection moves from E21 to G24.
#  - Sheet "Asset_Cal_Info": rows 16 and 18-24 have their Ref-Des/Mooring block
#    (columns A-D, plus G on rows 18-24) shifted right by 8 columns (A->I,
#    B->J, C->K, D->L, G->O) to make room for newly-populated OOI barcode /
#    NUTNR calibration columns, and the trailing "No calibration coefficient"
#    marker cell shifts from column I to column Q.

$wb  = $excel.ActiveWorkbook
$moorings = $wb.Worksheets.Item("Moorings")
$cal      = $wb.Worksheets.Item("Asset_Cal_Info")

# ---------------------------------------------------------------------------
# Helper: move (cut) the contents+formatting of $srcAddr to $dstAddr on sheet
# $ws, then fully clear the source cell (value, format, and number format) so
# that no stray empty cell element is left behind.
# ---------------------------------------------------------------------------
function Move-Cell {
    param($ws, [string]$srcAddr, [string]$dstAddr)
    $ws.Range($srcAddr).Cut($ws.Range($dstAddr))
    $ws.Range($srcAddr).Clear()
}

# ---------------------------------------------------------------------------
# 1) Moorings sheet: remove the lingering A2 selection by re-pointing the
#    view at the sheet's default cell (A1), then restoring Asset_Cal_Info as
#    the active/selected tab so the workbook's active-tab state is unchanged.
# ---------------------------------------------------------------------------
$moorings.Range("A1").Select()
$cal.Select()

# ---------------------------------------------------------------------------
# 2) Asset_Cal_Info: shift the block in row 16 eight columns to the right.
#    Move the existing Q-destination cell out of the way first so it is not
#    clobbered by the incoming data.
# ---------------------------------------------------------------------------
Move-Cell $cal "I16" "Q16"
Move-Cell $cal "A16" "I16"
Move-Cell $cal "B16" "J16"
Move-Cell $cal "C16" "K16"
Move-Cell $cal "D16" "L16"

# ---------------------------------------------------------------------------
# 3) Asset_Cal_Info: same eight-column shift for rows 18 through 24
#    (A->I, B->J, C->K, D->L, G->O, I->Q).
# ---------------------------------------------------------------------------
foreach ($r in 18..24) {
    Move-Cell $cal "I$r" "Q$r"
    Move-Cell $cal "A$r" "I$r"
    Move-Cell $cal "B$r" "J$r"
    Move-Cell $cal "C$r" "K$r"
    Move-Cell $cal "D$r" "L$r"
    Move-Cell $cal "G$r" "O$r"
}

# ---------------------------------------------------------------------------
# 4) Asset_Cal_Info: move the active selection from E21 to G24.
# ---------------------------------------------------------------------------
$cal.Range("G24").Select()
